$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "dolly1"
$ws.Range("B4").Value = "03ac674216f3e15c761ee1a5e255f067953623c8b388b4459e13f978d7c846f4"
$ws.Range("C4").Value = "dollydolly"
$ws.Range("D4").Value = "user"
$ws.Range("E4").Value = $true
